# envio semana 20 de 2025
# Adds a new "week 20" column (W) to the weekly IRA report: a header
# label "20" in W1 (matching the style of the existing week-number
# headers in row 1) plus the per-UPGD counts for that week in W2:W54
# (rows with no data for the week are intentionally left blank, mirroring
# the sparse V/U/T/... columns already in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell W1: text "20", same look as D1:V1 ("1".."19") ---
# Those headers are numeric-looking text, not numbers, so force the
# cell to Text before assigning, then re-pull the plain (unformatted)
# style from the neighboring header cell (V1) so W1 ends up using the
# exact same style as the rest of row 1 instead of a one-off "Text
# number format" style.
$ws.Cells.Item(1, 23).NumberFormat = "@"
$ws.Cells.Item(1, 23).Value = "20"
$ws.Cells.Item(1, 22).Copy()
$ws.Cells.Item(1, 23).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data cells W2:W54: week-20 counts per UPGD row ---
$weekValues = @{
    2  = 70
    4  = 0
    5  = 1
    6  = 90
    7  = 33
    8  = 37
    9  = 2
    10 = 1
    11 = 3
    12 = 4
    13 = 1
    14 = 3
    17 = 1
    20 = 1
    21 = 4
    22 = 1
    23 = 60
    24 = 3
    26 = 178
    27 = 0
    28 = 18
    29 = 3
    30 = 5
    32 = 26
    33 = 4
    34 = 12
    35 = 97
    36 = 1
    37 = 2
    38 = 49
    39 = 25
    41 = 87
    42 = 151
    43 = 6
    44 = 137
    45 = 1
    46 = 0
    47 = 5
    48 = 4
    49 = 54
    50 = 3
    51 = 0
    52 = 2
    53 = 17
    54 = 23
}

foreach ($row in $weekValues.Keys) {
    $ws.Cells.Item($row, 23).Value = $weekValues[$row]
}
